$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("metric overrides")
$ws6 = $wb.Worksheets.Item("production")

# ---------------------------------------------------------------------------
# sheet1 ("metric overrides"): the cohort-override table (rows 3-11) gets two
# new columns inserted before column D (old col B->D, C->E, D->F, E->G, F->H),
# and a brand-new row is inserted at row 7 ("Age months"), pushing the header
# + data rows (old 7-10) down to (8-11) and the blank placeholder rows below
# (old 11-21) down to (12-22), unshifted in columns. We replicate this purely
# by copying each source cell onto its destination cell (which carries value
# + number format/style together), walking bottom-right -> top-left so that
# every source is read before it is ever overwritten as someone else's dest.
# ---------------------------------------------------------------------------

# blank placeholder rows: old row N (B/D) -> new row N+1 (B/D), unshifted cols
$ws1.Range("D21").Copy($ws1.Range("D22"))
$ws1.Range("B21").Copy($ws1.Range("B22"))
$ws1.Range("D20").Copy($ws1.Range("D21"))
$ws1.Range("B20").Copy($ws1.Range("B21"))
$ws1.Range("D19").Copy($ws1.Range("D20"))
$ws1.Range("B19").Copy($ws1.Range("B20"))
$ws1.Range("D18").Copy($ws1.Range("D19"))
$ws1.Range("B18").Copy($ws1.Range("B19"))
$ws1.Range("D17").Copy($ws1.Range("D18"))
$ws1.Range("B17").Copy($ws1.Range("B18"))
$ws1.Range("D16").Copy($ws1.Range("D17"))
$ws1.Range("B16").Copy($ws1.Range("B17"))
$ws1.Range("D15").Copy($ws1.Range("D16"))
$ws1.Range("B15").Copy($ws1.Range("B16"))
$ws1.Range("D14").Copy($ws1.Range("D15"))
$ws1.Range("B14").Copy($ws1.Range("B15"))
$ws1.Range("D13").Copy($ws1.Range("D14"))
$ws1.Range("B13").Copy($ws1.Range("B14"))
$ws1.Range("D12").Copy($ws1.Range("D13"))
$ws1.Range("B12").Copy($ws1.Range("B13"))
$ws1.Range("D11").Copy($ws1.Range("D12"))
$ws1.Range("B11").Copy($ws1.Range("B12"))

# data row: old row 10 -> new row 11 (cols B-F -> D-H)
$ws1.Range("F10").Copy($ws1.Range("H11"))
$ws1.Range("E10").Copy($ws1.Range("G11"))
$ws1.Range("D10").Copy($ws1.Range("F11"))
$ws1.Range("C10").Copy($ws1.Range("E11"))
$ws1.Range("B10").Copy($ws1.Range("D11"))
$ws1.Range("A10").Copy($ws1.Range("A11"))

# data row: old row 9 -> new row 10 (cols B-F -> D-H)
$ws1.Range("F9").Copy($ws1.Range("H10"))
$ws1.Range("E9").Copy($ws1.Range("G10"))
$ws1.Range("D9").Copy($ws1.Range("F10"))
$ws1.Range("C9").Copy($ws1.Range("E10"))
$ws1.Range("B9").Copy($ws1.Range("D10"))
$ws1.Range("A9").Copy($ws1.Range("A10"))

# data row: old row 8 -> new row 9 (cols B-F -> D-H)
$ws1.Range("F8").Copy($ws1.Range("H9"))
$ws1.Range("E8").Copy($ws1.Range("G9"))
$ws1.Range("D8").Copy($ws1.Range("F9"))
$ws1.Range("C8").Copy($ws1.Range("E9"))
$ws1.Range("B8").Copy($ws1.Range("D9"))
$ws1.Range("A8").Copy($ws1.Range("A9"))

# header row: old row 7 -> new row 8 (cols B-F -> D-H)
$ws1.Range("F7").Copy($ws1.Range("H8"))
$ws1.Range("E7").Copy($ws1.Range("G8"))
$ws1.Range("D7").Copy($ws1.Range("F8"))
$ws1.Range("C7").Copy($ws1.Range("E8"))
$ws1.Range("B7").Copy($ws1.Range("D8"))
$ws1.Range("A7").Copy($ws1.Range("A8"))

# rows 3-6 (no row shift, only the col B:C insert): cols B-F -> D-H
$ws1.Range("F6").Copy($ws1.Range("H6"))
$ws1.Range("B6").Copy($ws1.Range("D6"))
$ws1.Range("F5").Copy($ws1.Range("H5"))
$ws1.Range("E5").Copy($ws1.Range("G5"))
$ws1.Range("D5").Copy($ws1.Range("F5"))
$ws1.Range("C5").Copy($ws1.Range("E5"))
$ws1.Range("B5").Copy($ws1.Range("D5"))
$ws1.Range("F4").Copy($ws1.Range("H4"))
$ws1.Range("E4").Copy($ws1.Range("G4"))
$ws1.Range("D4").Copy($ws1.Range("F4"))
$ws1.Range("C4").Copy($ws1.Range("E4"))
$ws1.Range("B4").Copy($ws1.Range("D4"))
$ws1.Range("F3").Copy($ws1.Range("H3"))
$ws1.Range("E3").Copy($ws1.Range("G3"))
$ws1.Range("D3").Copy($ws1.Range("F3"))
$ws1.Range("C3").Copy($ws1.Range("E3"))
$ws1.Range("B3").Copy($ws1.Range("D3"))

# wipe the stale leftovers of the old layout (content AND formatting - these
# cells must not exist at all afterwards) that are not reused as a
# destination above and are not about to be overwritten with new content below
foreach ($addr in @("B3","C3","B4","C4","B5","C5","B6","F6","A7","B7","C7","E7","F7","B9","C9","B10","C10","B11")) {
    $ws1.Range($addr).Clear()
}

# new row 7: "Age months" cohort-filter header label
$ws1.Range("D7").Value = "Age months"

# new header cells for the inserted columns
$ws1.Range("B8").Value = "Book"
$ws1.Range("C8").Value = "MinAgeMonths"

# new cohort condition on the (now) row 10 data row: only apply to items >= 36 months old
$ws1.Range("C10").Value = 36

# new cohort row (row 11): a "new" book override, effective 2025-10-15
$ws1.Range("B11").Value = "new"
$ws1.Range("D11").Value = 45945

# ---------------------------------------------------------------------------
# sheet6 ("production"): bump the first mutation's date and append a new one
# ---------------------------------------------------------------------------
$ws6.Range("A3").Value = 45930

$ws6.Range("A3:D3").Copy($ws6.Range("A6:D6"))
$ws6.Range("A6").Value = 46356
$ws6.Range("C6").Value = 20000

# ---------------------------------------------------------------------------
# selection / active sheet bookkeeping to mirror the saved view state
# ---------------------------------------------------------------------------
$ws1.Range("C11").Select() | Out-Null
$ws6.Range("A4").Select() | Out-Null
$ws6.Activate() | Out-Null
